$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19: Inscritos 31->32, Pagos 19->20, Inscrições homologadas 19->20
$ws.Range("E19").Value = 32
$ws.Range("F19").Value = 20
$ws.Range("H19").Value = 20

# Row 28: Inscritos 7->8
$ws.Range("E28").Value = 8

# Row 55: Inscritos 3->4
$ws.Range("E55").Value = 4

# Row 65: Inscritos 22->23
$ws.Range("E65").Value = 23

# Row 66: Inscritos 25->26
$ws.Range("E66").Value = 26

# Row 87: Inscritos 5->6
$ws.Range("E87").Value = 6
